$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 42.409254
$ws.Range("H2").Value = 127.227762
$ws.Range("I2").Value = 0.6138221220752584
$ws.Range("J2").Value = 0.6138221220752584
$ws.Range("M2").Value = 29.52617166666667
$ws.Range("N2").Value = 88.57851500000001
$ws.Range("O2").Value = 0.3218391660320701
$ws.Range("P2").Value = 0.3218391660320701
$ws.Range("Q2").Value = 1252.18291385927
$ws.Range("R2").Value = 11269.64622473343
$ws.Range("S2").Value = 0.1975519998607367
$ws.Range("T2").Value = 0.1975519998607367
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 42.409254
$ws.Range("H3").Value = 127.227762
$ws.Range("I3").Value = 0.6138221220752584
$ws.Range("J3").Value = 0.6138221220752584
$ws.Range("O3").Value = 0.4328989896002822
$ws.Range("P3").Value = 0.4328989896002822
$ws.Range("Q3").Value = 1684.284498022842
$ws.Range("R3").Value = 15158.56048220558
$ws.Range("S3").Value = 0.2657229764406805
$ws.Range("T3").Value = 0.2657229764406805
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 42.409254
$ws.Range("H4").Value = 127.227762
$ws.Range("I4").Value = 0.6138221220752584
$ws.Range("J4").Value = 0.6138221220752584
$ws.Range("M4").Value = 22.50081433333333
$ws.Range("N4").Value = 67.502443
$ws.Range("O4").Value = 0.2452618443676477
$ws.Range("P4").Value = 0.2452618443676476
$ws.Range("Q4").Value = 954.2427502691742
$ws.Range("R4").Value = 8588.184752422567
$ws.Range("S4").Value = 0.1505471457738413
$ws.Range("T4").Value = 0.1505471457738412
$ws.Range("I5").Value = 0.07014398987036251
$ws.Range("J5").Value = 0.07014398987036251
$ws.Range("M5").Value = 29.52617166666667
$ws.Range("N5").Value = 88.57851500000001
$ws.Range("O5").Value = 0.3218391660320701
$ws.Range("P5").Value = 0.3218391660320701
$ws.Range("Q5").Value = 143.0921149088478
$ws.Range("R5").Value = 1287.82903417963
$ws.Range("S5").Value = 0.02257508320203944
$ws.Range("T5").Value = 0.02257508320203944
$ws.Range("I6").Value = 0.07014398987036251
$ws.Range("J6").Value = 0.07014398987036251
$ws.Range("O6").Value = 0.4328989896002822
$ws.Range("P6").Value = 0.4328989896002822
$ws.Range("S6").Value = 0.03036526234141236
$ws.Range("T6").Value = 0.03036526234141236
$ws.Range("I7").Value = 0.07014398987036251
$ws.Range("J7").Value = 0.07014398987036251
$ws.Range("M7").Value = 22.50081433333333
$ws.Range("N7").Value = 67.502443
$ws.Range("O7").Value = 0.2452618443676477
$ws.Range("P7").Value = 0.2452618443676476
$ws.Range("Q7").Value = 109.0452614878895
$ws.Range("R7").Value = 981.4073533910059
$ws.Range("S7").Value = 0.0172036443269107
$ws.Range("T7").Value = 0.0172036443269107
$ws.Range("G8").Value = 21.83492733333334
$ws.Range("H8").Value = 65.50478200000001
$ws.Range("I8").Value = 0.3160338880543792
$ws.Range("J8").Value = 0.3160338880543791
$ws.Range("M8").Value = 29.52617166666667
$ws.Range("N8").Value = 88.57851500000001
$ws.Range("O8").Value = 0.3218391660320701
$ws.Range("P8").Value = 0.3218391660320701
$ws.Range("Q8").Value = 644.7018127731924
$ws.Range("R8").Value = 5802.316314958731
$ws.Range("S8").Value = 0.101712082969294
$ws.Range("T8").Value = 0.101712082969294
$ws.Range("G9").Value = 21.83492733333334
$ws.Range("H9").Value = 65.50478200000001
$ws.Range("I9").Value = 0.3160338880543792
$ws.Range("J9").Value = 0.3160338880543791
$ws.Range("O9").Value = 0.4328989896002822
$ws.Range("P9").Value = 0.4328989896002822
$ws.Range("Q9").Value = 867.1746412466622
$ws.Range("R9").Value = 7804.571771219959
$ws.Range("S9").Value = 0.1368107508181894
$ws.Range("T9").Value = 0.1368107508181894
$ws.Range("G10").Value = 21.83492733333334
$ws.Range("H10").Value = 65.50478200000001
$ws.Range("I10").Value = 0.3160338880543792
$ws.Range("J10").Value = 0.3160338880543791
$ws.Range("M10").Value = 22.50081433333333
$ws.Range("N10").Value = 67.502443
$ws.Range("O10").Value = 0.2452618443676477
$ws.Range("P10").Value = 0.2452618443676476
$ws.Range("Q10").Value = 491.3036459091585
$ws.Range("R10").Value = 4421.732813182426
$ws.Range("S10").Value = 0.07751105426689572
$ws.Range("T10").Value = 0.07751105426689571
